# Update gh-pages to output generated at 456a3b4
# Applies the same set of data refreshes to both the "展览" sheet and the
# "全部类型" sheet (which duplicates the first 16 data rows of "展览").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 (ACGN summer pool party): "想去人数" (interested count) 625 -> 626
    $ws.Range("F2").Value = 626

    # Row 3 (IE动漫嘉年华): 想去人数 575 -> 576, and cover image URL updated
    $ws.Range("F3").Value = 576
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202409/DoyjJqXA1725871233303.jpeg"

    # Row 6 (银泰百货高新店嘉年华): 想去人数 119 -> 120
    $ws.Range("F6").Value = 120

    # Row 10 (第十五届次元之门动漫游戏博览会): 想去人数 5014 -> 5030, 最低票价 70 -> 59.9
    $ws.Range("F10").Value = 5030
    $ws.Range("G10").Value = 59.9

    # Row 11 (首届AT次元时代动漫游戏嘉年华): 想去人数 4701 -> 4707
    $ws.Range("F11").Value = 4707

    # Row 16 (W·A第五人格同人only2.0): 想去人数 174 -> 175
    $ws.Range("F16").Value = 175
}
